$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) share identical rows 2-16
# and both need the same "想去人数" (column F) counts updated.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 1167
    $ws.Range("F3").Value = 590
    $ws.Range("F7").Value = 61
    $ws.Range("F10").Value = 5392
    $ws.Range("F11").Value = 4847
}
